# Generate Report for Handback
#
# Overview sheet: "Status" column (E/F) flips from "Ready for handoff" to
# "Handed back: in sync with en-US". The same status text is shown on the
# zh-cn / de-de detail sheets (column C).
#
# zh-cn / de-de detail sheets: the handback run now has a Target File
# (hyperlinked .md), a Handback File (the generated .xlf) and a Handback
# DateTime stamped for each locale.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$mdFile  = "830954f8-9fd3-49ce-96f5-15d7d3de2793.md"
$mdUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/516380630c823f867d76ad23d9965b85a857070f/e2e/830954f8-9fd3-49ce-96f5-15d7d3de2793.md"
$zhXlf   = "830954f8-9fd3-49ce-96f5-15d7d3de2793.d1633fd436dd871822e2d1fa9bebe567e6e2583d.zh-cn.xlf"
$deXlf   = "830954f8-9fd3-49ce-96f5-15d7d3de2793.d1633fd436dd871822e2d1fa9bebe567e6e2583d.de-de.xlf"
$statusText = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # RGB(100,149,237) == FF6495ED, matches the workbook's HyperLink style

# ---------------------------------------------------------------------
# Overview: Status text for both locale columns
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = $statusText
$ws1.Range("F2").Value = $statusText
$ws1.Columns.Item(5).ColumnWidth = 29.15
$ws1.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$ws2.Range("C2").Value = $statusText

$ws2.Range("I2").Value = $mdFile
$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdUrl, "", "", $mdFile)
$ws2.Range("I2").Font.Name = "Calibri"
$ws2.Range("I2").Font.Size = 11
$ws2.Range("I2").Font.Underline = 2
$ws2.Range("I2").Font.Color = $hyperlinkColor

$ws2.Range("J2").Value = $zhXlf
$ws2.Range("K2").Value = "2016-08-18 21:00:07"

$ws2.Columns.Item(3).ColumnWidth = 29.15
$ws2.Columns.Item(9).ColumnWidth = 39.15
$ws2.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$ws3.Range("C2").Value = $statusText

$ws3.Range("I2").Value = $mdFile
$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdUrl, "", "", $mdFile)
$ws3.Range("I2").Font.Name = "Calibri"
$ws3.Range("I2").Font.Size = 11
$ws3.Range("I2").Font.Underline = 2
$ws3.Range("I2").Font.Color = $hyperlinkColor

$ws3.Range("J2").Value = $deXlf
$ws3.Range("K2").Value = "2016-08-18 21:00:29"

$ws3.Columns.Item(3).ColumnWidth = 29.15
$ws3.Columns.Item(9).ColumnWidth = 39.15
$ws3.Columns.Item(10).ColumnWidth = 39.15
